$d = $word.ActiveDocument

# Title paragraph: merge the word-by-word runs into a single run
$d.Content.Find.Execute(
    "Answers: Introduction to sigma notation", $true, $false, $false, $false, $false,
    $true, 1, $false, "Answers: Introduction to sigma notation", 2) | Out-Null

# Author paragraph: merge the word-by-word runs into a single run
$d.Content.Find.Execute(
    "Ifan Howells-Baines, Mark Toner", $true, $false, $false, $false, $false,
    $true, 1, $false, "Ifan Howells-Baines, Mark Toner", 2) | Out-Null

# Abstract paragraph: merge the word-by-word runs into a single run
$d.Content.Find.Execute(
    "Answers to questions relating to the guide on introduction to sigma notation.", $true, $false, $false, $false, $false,
    $true, 1, $false, "Answers to questions relating to the guide on introduction to sigma notation.", 2) | Out-Null
